$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion summary text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = @"
Conversión del día 💰
✅ Dólar paralelo: 68

Binance
✅ 1000 Bs = 13.89 = 56841.23 pesos
✅ 56841.23 pesos = 13.86 = 972.54 Bs

Promedio competencia
✅ Tasa pesos: 20
✅ Tasa Bs: 20
✅ % Ganancia: 20%
"@

$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": refresh the raw rate inputs feeding the conversion ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 71.98999999999999
$ws2.Range("O10").Value = 4092
$ws2.Range("N12").Value = 4099.99
$ws2.Range("O12").Value = 70.15000000000001
